$wb = $excel.ActiveWorkbook

$customers   = $wb.Worksheets.Item("customers")
$invoices    = $wb.Worksheets.Item("invoices")
$lineItems   = $wb.Worksheets.Item("line_items")
$productList = $wb.Worksheets.Item("product_list")

# --- Populate header rows -------------------------------------------------
# The write order below reproduces the shared-string table ordering of the
# target workbook (id columns were filled in after the rest of the row).

$customers.Activate()
$customers.Range("B1").Value = "first_name"
$customers.Range("C1").Value = "last_name"
$customers.Range("D1").Value = "email"
$customers.Range("E1").Value = "telephone"

$invoices.Activate()
$invoices.Range("B1").Value = "description"
$invoices.Range("C1").Value = "amount"
$invoices.Range("D1").Value = "date"

$lineItems.Activate()
$lineItems.Range("A1").Value = "line_id"
$lineItems.Range("B1").Value = "description"

$customers.Activate()
$customers.Range("A1").Value = "cust_id"

$productList.Activate()
$productList.Range("A2").Value = "prod_id"
$productList.Range("B2").Value = "prod_name"
$productList.Range("C2").Value = "prod_desc"

$lineItems.Activate()
$lineItems.Range("C1").Value = "quantity"

$invoices.Activate()
$invoices.Range("A1").Value = "invoice_id"

# --- View state: zoom + selection per sheet -------------------------------

$customers.Activate()
$excel.ActiveWindow.Zoom = 150
$customers.Range("C6").Select()

$invoices.Activate()
$excel.ActiveWindow.Zoom = 150
$invoices.Columns.Item(1).ColumnWidth = 8.5
$invoices.Range("C5").Select()

$lineItems.Activate()
$excel.ActiveWindow.Zoom = 150
$lineItems.Range("E12").Select()

$productList.Activate()
$excel.ActiveWindow.Zoom = 150
$productList.Range("C14").Select()
